$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" '29.879.54'
Set-TextValue "E2" '  +0.46%  '

Set-TextValue "D3" '1.894.97'
Set-TextValue "E3" '  +0.46%  '

Set-TextValue "E4" '  +0.14%  '

Set-TextValue "D5" '0.7826'
Set-TextValue "E5" '  -1.34%  '

Set-TextValue "D6" '243.76'
Set-TextValue "E6" '  +1.02%  '

Set-TextValue "D8" '0.3137'
Set-TextValue "E8" '  -1.08%  '

Set-TextValue "D9" '25.68'
Set-TextValue "E9" '  +0.94%  '

Set-TextValue "D10" '0.07273'
Set-TextValue "E10" '  +3.86%  '

Set-TextValue "D11" '0.08106'
Set-TextValue "E11" '  +0.91%  '

Set-TextValue "D12" '0.7734'
Set-TextValue "E12" '  +1.55%  '

Set-TextValue "D13" '5.472'
Set-TextValue "E13" '  +3.31%  '

Set-TextValue "D14" '1.888.56'
Set-TextValue "E14" '  -0.77%  '

Set-TextValue "D15" '94.66'
Set-TextValue "E15" '  +2.64%  '

Set-TextValue "D16" '6.207'
Set-TextValue "E16" '  +4.66%  '

Set-TextValue "D17" '29.872.70'
Set-TextValue "E17" '  +0.59%  '

Set-TextValue "E18" '  +0.78%  '

Set-TextValue "D19" '246.02'
Set-TextValue "E19" '  +1.09%  '

Set-TextValue "D20" '0.000007834'
Set-TextValue "E20" '  +2.04%  '

Set-TextValue "E21" '  +0.11%  '

Set-TextValue "D22" '8.135'
Set-TextValue "E22" '  -1.03%  '

Set-TextValue "D23" '2.136.25'
Set-TextValue "E23" '  +2.60%  '

Set-TextValue "D24" '1.002'
Set-TextValue "E24" '  +0.17%  '

Set-TextValue "D25" '0.1593'
Set-TextValue "E25" '  -5.40%  '

Set-TextValue "D26" '9.459'
Set-TextValue "E26" '  +1.82%  '

Set-TextValue "D27" '164.07'
Set-TextValue "E27" '  -0.08%  '

Set-TextValue "D28" '18.78'
Set-TextValue "E28" '  +0.91%  '

Set-TextValue "D29" '2.021'
Set-TextValue "E29" '  -1.36%  '

Set-TextValue "D30" '1.436'

Set-TextValue "D31" '1.543'
Set-TextValue "E31" '  +0.74%  '

Set-TextValue "D32" '4.471'
Set-TextValue "E32" '  +2.18%  '

Set-TextValue "D33" '0.05573'
Set-TextValue "E33" '  -1.78%  '

Set-TextValue "D34" '4.076'
Set-TextValue "E34" '  +0.67%  '

Set-TextValue "D35" '1.244'
Set-TextValue "E35" '  -1.34%  '

Set-TextValue "D36" '0.7530'
Set-TextValue "E36" '  +2.60%  '

Set-TextValue "D37" '1.002'
Set-TextValue "E37" '  +0.70%  '

Set-TextValue "D38" '2.680'
Set-TextValue "E38" '  +2.48%  '

Set-TextValue "D39" '0.01934'
Set-TextValue "E39" '  +1.39%  '

Set-TextValue "E40" '  +0.67%  '

Set-TextValue "D41" '1.143.70'
Set-TextValue "E41" '  +12.11%  '

Set-TextValue "D42" '0.4457'
Set-TextValue "E42" '  +1.22%  '

Set-TextValue "D43" '74.00'
Set-TextValue "E43" '  +2.18%  '

Set-TextValue "D44" '5.960'
Set-TextValue "E44" '  +2.43%  '

Set-TextValue "D45" '0.8526'
Set-TextValue "E45" '  +2.11%  '

Set-TextValue "E46" '  +0.15%  '

Set-TextValue "D47" '1.895'
Set-TextValue "E47" '  +1.64%  '

Set-TextValue "D48" '3.156'
Set-TextValue "E48" '  +8.77%  '

Set-TextValue "D49" '102.07'
Set-TextValue "E49" '  -0.49%  '

Set-TextValue "E50" '  +1.58%  '

Set-TextValue "D51" '9.767'
Set-TextValue "E51" '  -1.01%  '
